$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 479.33334
$ws.Range("J12").Value = 219.5
$ws.Range("L12").Value = 219.5
$ws.Range("N12").Value = -559.5
$ws.Range("H21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()
$ws.Range("H23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()
$ws.Range("H29").Value = 3995.4
$ws.Range("J29").Value = 3995.4
$ws.Range("L29").Value = 11986.2
$ws.Range("N29").Value = -12548.2
$ws.Range("H38").Value = 2566.6667
$ws.Range("J38").Value = 6500
$ws.Range("L38").Value = 19500
$ws.Range("N38").Value = -20244
$ws.Range("H41").Value = 649.25
$ws.Range("J41").Value = 500
$ws.Range("L41").Value = 500
$ws.Range("N41").Value = -1380
$ws.Range("H58").Value = 7673.7144
$ws.Range("J58").Value = 10583.2
$ws.Range("L58").Value = 31749.6
$ws.Range("N58").Value = -32049.6
$ws.Range("H61").Value = 694
$ws.Range("I61").Value = 694
$ws.Range("K61").Value = 2082
$ws.Range("M61").Value = -1910
$ws.Range("H137").Value = 2296.8125
$ws.Range("I137").Value = 1635.2858
$ws.Range("K137").Value = 4905.857400000001
$ws.Range("M137").Value = -2355.857400000001
$ws.Range("H138").Value = 3352
$ws.Range("J138").Value = 3362.3
$ws.Range("L138").Value = 10086.9
$ws.Range("N138").Value = -20366.9

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1183.5714
$ws.Range("I2").Value = 1113.9375
$ws.Range("J2").Value = 1406.4
$ws.Range("K2").Value = 1113.9375
$ws.Range("L2").Value = 1406.4
$ws.Range("M2").Value = -1000.9375
$ws.Range("N2").Value = -1632.4
$ws.Range("H32").Value = 7272.5415
$ws.Range("I32").Value = 5099.6343
$ws.Range("J32").Value = 19999.572
$ws.Range("K32").Value = 5099.6343
$ws.Range("L32").Value = 19999.572
$ws.Range("M32").Value = -4812.6343
$ws.Range("N32").Value = -20573.572
$ws.Range("H74").Value = 800
$ws.Range("I74").Value = 800
$ws.Range("K74").Value = 800
$ws.Range("M74").Value = 74
$ws.Range("H77").Value = 800
$ws.Range("I77").Value = 800
$ws.Range("K77").Value = 4000
$ws.Range("M77").Value = 368
$ws.Range("H97").Value = 219.46153
$ws.Range("I97").Value = 154.41667
$ws.Range("K97").Value = 154.41667
$ws.Range("M97").Value = 341.58333
$ws.Range("H116").Value = 1183.5714
$ws.Range("I116").Value = 1113.9375
$ws.Range("J116").Value = 1406.4
$ws.Range("K116").Value = 1113.9375
$ws.Range("L116").Value = 1406.4
$ws.Range("M116").Value = 1180.0625
$ws.Range("N116").Value = -5994.4

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1183.5714
$ws.Range("I3").Value = 1113.9375
$ws.Range("J3").Value = 1406.4
$ws.Range("K3").Value = 1113.9375
$ws.Range("L3").Value = 1406.4
$ws.Range("M3").Value = -999.9375
$ws.Range("N3").Value = -1634.4
$ws.Range("H86").Value = 2547.7273
$ws.Range("J86").Value = 2675
$ws.Range("L86").Value = 2675
$ws.Range("N86").Value = -4921
$ws.Range("H89").Value = 2547.7273
$ws.Range("J89").Value = 2675
$ws.Range("L89").Value = 13375
$ws.Range("N89").Value = -24607
$ws.Range("H94").Value = 1027.5834
$ws.Range("I94").Value = 566.5
$ws.Range("K94").Value = 566.5
$ws.Range("M94").Value = -115.5
$ws.Range("H105").Value = 5672
$ws.Range("J105").Value = 6722
$ws.Range("L105").Value = 6722
$ws.Range("N105").Value = -10216
$ws.Range("H107").Value = 377.66666
$ws.Range("I107").Value = 331.7143
$ws.Range("J107").Value = 538.5
$ws.Range("K107").Value = 331.7143
$ws.Range("L107").Value = 538.5
$ws.Range("M107").Value = 1588.2857
$ws.Range("N107").Value = -4378.5

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1599
$ws.Range("I31").Value = 1599.5
$ws.Range("J31").Value = 1598
$ws.Range("K31").Value = 1599.5
$ws.Range("L31").Value = 1598
$ws.Range("M31").Value = -1304.5
$ws.Range("N31").Value = -2188
$ws.Range("H34").Value = 1599
$ws.Range("I34").Value = 1599.5
$ws.Range("J34").Value = 1598
$ws.Range("K34").Value = 1599.5
$ws.Range("L34").Value = 1598
$ws.Range("M34").Value = -1397.5
$ws.Range("N34").Value = -2002
$ws.Range("H99").Value = 4489
$ws.Range("I99").Value = 4316.5
$ws.Range("J99").Value = 5006.5
$ws.Range("K99").Value = 4316.5
$ws.Range("L99").Value = 5006.5
$ws.Range("M99").Value = -2818.5
$ws.Range("N99").Value = -8002.5
$ws.Range("H122").Value = 4499.75
$ws.Range("I122").Value = 3999.5
$ws.Range("K122").Value = 11998.5
$ws.Range("M122").Value = -9548.5
$ws.Range("H126").Value = 4489
$ws.Range("I126").Value = 4316.5
$ws.Range("J126").Value = 5006.5
$ws.Range("K126").Value = 12949.5
$ws.Range("L126").Value = 15019.5
$ws.Range("M126").Value = -10479.5
$ws.Range("N126").Value = -19959.5
$ws.Range("H132").Value = 3858.182
$ws.Range("I132").Value = 3332.6667
$ws.Range("K132").Value = 9998.000100000001
$ws.Range("M132").Value = -7468.000100000001
$ws.Range("H134").Value = 3458.6667
$ws.Range("I134").Value = 3458.6667
$ws.Range("K134").Value = 10376.0001
$ws.Range("M134").Value = -7841.000100000001

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 146998.58
$ws.Range("J80").Value = 170665.67
$ws.Range("L80").Value = 511997.01
$ws.Range("N80").Value = -513869.01
$ws.Range("H83").Value = 146998.58
$ws.Range("J83").Value = 170665.67
$ws.Range("L83").Value = 1535991.03
$ws.Range("N83").Value = -1545351.03

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5998.3335
$ws.Range("I80").Value = 3994.5
$ws.Range("K80").Value = 3994.5
$ws.Range("M80").Value = -2996.5
$ws.Range("H83").Value = 5998.3335
$ws.Range("I83").Value = 3994.5
$ws.Range("K83").Value = 19972.5
$ws.Range("M83").Value = -14980.5
$ws.Range("H102").Value = 1370.3846
$ws.Range("J102").Value = 2033
$ws.Range("L102").Value = 2033
$ws.Range("N102").Value = -5277
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("H132").Value = 3776
$ws.Range("I132").Value = 3698.3333
$ws.Range("K132").Value = 11094.9999
$ws.Range("M132").Value = -8564.999899999999

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("H132").Value = 4979.25
$ws.Range("I132").Value = 3959.5
$ws.Range("K132").Value = 11878.5
$ws.Range("M132").Value = -9348.5
$ws.Range("H136").Value = 4677.7144
$ws.Range("I136").Value = 5124.3335
$ws.Range("J136").Value = 1998
$ws.Range("K136").Value = 15373.0005
$ws.Range("L136").Value = 5994
$ws.Range("M136").Value = -12823.0005
$ws.Range("N136").Value = -11094

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1499
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 1499
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 4497
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -9437
$ws.Range("H132").Value = 2504.2307
$ws.Range("I132").Value = 2283.889
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 6851.667
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -4321.667
$ws.Range("N132").Value = -14060
$ws.Range("H136").Value = 2409.6
$ws.Range("I136").Value = 1969.45
$ws.Range("K136").Value = 5908.35
$ws.Range("M136").Value = -3358.35
$ws.Range("H139").Value = 69800
$ws.Range("J139").Value = 69800
$ws.Range("L139").Value = 69800
$ws.Range("N139").Value = -80080
